# Atualização de bases das ligas, do dia: 11-06-2024 às 21:19
# Swap the betting-odds data between rows that belong to the same fixture
# (Div/Date match) but whose HomeTeam/AwayTeam/stat columns got mixed up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($ws, $row1, $row2, $firstCol, $lastCol) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Columns B (2) through AD (30) are swapped; A (id) and C/D (Div/Date) stay put.
$firstCol = 2   # B
$lastCol  = 30  # AD

# Rows 117 <-> 118
Swap-RowRange $ws 117 118 $firstCol $lastCol

# Rows 179 <-> 180
Swap-RowRange $ws 179 180 $firstCol $lastCol

# Rows 234, 235, 237 rotate: 234 <- 235, 235 <- 237, 237 <- 234
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $v234 = $ws.Cells.Item(234, $c).Value()
    $v235 = $ws.Cells.Item(235, $c).Value()
    $v237 = $ws.Cells.Item(237, $c).Value()

    $ws.Cells.Item(234, $c).Value = $v235
    $ws.Cells.Item(235, $c).Value = $v237
    $ws.Cells.Item(237, $c).Value = $v234
}
